# C5-PowerPoint.pptx — Fri, May 08, 2020  4:06:23 AM
#
# 1) Slide 6's table switches from the deck's custom "Table_0" style to the
#    built-in PowerPoint table style {D494E503-E332-4CB8-91E3-98303F80F036}.
# 2) The slide master's theme is repainted from the "Integral" colour
#    scheme to the standard "Office" colour scheme (dk1/lt1/dk2/lt2/accent1-6/
#    hlink/folHlink), i.e. the deck now renders with Office Theme colours.

$p = $ppt.ActivePresentation

# --- 1) Table style -------------------------------------------------------
$tableSlide = $p.Slides.Item(6)
$tableShape = $tableSlide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{D494E503-E332-4CB8-91E3-98303F80F036}")

# --- 2) Theme colours: Integral -> Office ---------------------------------
function HexToVbaRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $b * 65536 + $g * 256 + $r
}

# Order matches the 12 ThemeColorScheme slots: dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink.
$officeThemeColors = @(
    "000000", "FFFFFF", "44546A", "E7E6E6",
    "5B9BD5", "ED7D31", "A5A5A5", "FFC000",
    "4472C4", "70AD47", "0563C1", "954F72"
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Colors($i).RGB = HexToVbaRgb($officeThemeColors[$i - 1])
}
